$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy formatting from column R (2021) into the new column S (2022)
#    so the new cells inherit the same number formats / fonts / borders.
$ws.Range("R4:R34").Copy()
$ws.Range("S4:S34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 24 used a slightly different (but visually equivalent) style on R;
# the new S24 cell uses the more common style instead, matching R23.
$ws.Range("R23").Copy()
$ws.Range("S24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2. Fill in the new column S values (2022 data), row by row.
$ws.Range("S4").Value = 2022

$ws.Range("S5").Value = 0.5
$ws.Range("S6").Value = 0.2
$ws.Range("S7").Value = 0.7
$ws.Range("S8").Value = 0.2
$ws.Range("S9").Value = "-"
$ws.Range("S10").Value = 0.4
$ws.Range("S11").Value = 0.5
$ws.Range("S12").Value = 0.3
$ws.Range("S13").Value = 0.6
$ws.Range("S14").Value = 0.7
$ws.Range("S15").Value = 0.4
$ws.Range("S16").Value = 1.1000000000000001
$ws.Range("S17").Value = "-"
$ws.Range("S18").Value = "-"
$ws.Range("S19").Value = "-"
$ws.Range("S20").Value = 0.4
$ws.Range("S21").Value = 0.4
$ws.Range("S22").Value = 0.4
$ws.Range("S23").Value = 0.4
$ws.Range("S24").Value = "-"
$ws.Range("S25").Value = 0.7
$ws.Range("S26").Value = 1
$ws.Range("S27").Value = 0.4
$ws.Range("S28").Value = 1.7
$ws.Range("S29").Value = 0.3
$ws.Range("S30").Value = 0
$ws.Range("S31").Value = 0.6
$ws.Range("S32").Value = "-"
$ws.Range("S33").Value = "-"
$ws.Range("S34").Value = "-"

# 3. Move the active selection to T6, matching the post-edit workbook state.
$ws.Range("T6").Select()
